$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with new columns O1:R1 (copy style from N1, then set values)
$ws.Range("N1").Copy($ws.Range("O1:R1"))
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16

# Update/add data values for rows 2-25
$ws.Range("C2").Value = 1.012427397163468
$ws.Range("D2").Value = 1.029580967166665
$ws.Range("E2").Value = 1.017197308050836
$ws.Range("I2").Value = 1.050860526309924
$ws.Range("J2").Value = 1.034242655659945
$ws.Range("K2").Value = 1.040643018470422
$ws.Range("L2").Value = 1.028422145343317
$ws.Range("N2").Value = 1.015288783943612
$ws.Range("Q2").Value = 1.02
$ws.Range("R2").Value = 1.03980841480658
$ws.Range("C3").Value = 1.016338439300783
$ws.Range("D3").Value = 1.032198620421329
$ws.Range("E3").Value = 1.020284951199095
$ws.Range("I3").Value = 1.051998408056016
$ws.Range("J3").Value = 1.036393399089924
$ws.Range("K3").Value = 1.042435445253733
$ws.Range("L3").Value = 1.030664698014449
$ws.Range("N3").Value = 1.016005437727368
$ws.Range("Q3").Value = 1.02
$ws.Range("R3").Value = 1.041073152602141
$ws.Range("C4").Value = 1.018822179522997
$ws.Range("D4").Value = 1.033864384777591
$ws.Range("E4").Value = 1.022251514340127
$ws.Range("I4").Value = 1.052711208085103
$ws.Range("J4").Value = 1.037756366251449
$ws.Range("K4").Value = 1.043570579492924
$ws.Range("L4").Value = 1.03208860062446
$ws.Range("N4").Value = 1.01645980665881
$ws.Range("Q4").Value = 1.02
$ws.Range("R4").Value = 1.041876680552716
$ws.Range("C5").Value = 1.019859152926938
$ws.Range("D5").Value = 1.034562766834991
$ws.Range("E5").Value = 1.023074385939267
$ws.Range("I5").Value = 1.053008020148815
$ws.Range("J5").Value = 1.038326279776767
$ws.Range("K5").Value = 1.044046600270033
$ws.Range("L5").Value = 1.032684252174707
$ws.Range("N5").Value = 1.016650430767136
$ws.Range("Q5").Value = 1.02
$ws.Range("R5").Value = 1.042220414986112
$ws.Range("C6").Value = 1.020037142743111
$ws.Range("D6").Value = 1.034685276764056
$ws.Range("E6").Value = 1.023216282334844
$ws.Range("I6").Value = 1.053060780878618
$ws.Range("J6").Value = 1.038425984401868
$ws.Range("K6").Value = 1.04413175831098
$ws.Range("L6").Value = 1.032788006094072
$ws.Range("N6").Value = 1.016684490278008
$ws.Range("Q6").Value = 1.02
$ws.Range("R6").Value = 1.04228934656243
$ws.Range("C7").Value = 1.018848166604012
$ws.Range("D7").Value = 1.033888808196684
$ws.Range("E7").Value = 1.022273691637428
$ws.Range("I7").Value = 1.052723884076417
$ws.Range("J7").Value = 1.037775802196327
$ws.Range("K7").Value = 1.043591861324034
$ws.Range("L7").Value = 1.032107610669394
$ws.Range("N7").Value = 1.016468229035537
$ws.Range("Q7").Value = 1.02
$ws.Range("R7").Value = 1.04191178262118
$ws.Range("C8").Value = 1.013774042734003
$ws.Range("D8").Value = 1.030490058077929
$ws.Range("E8").Value = 1.018261258383744
$ws.Range("I8").Value = 1.051260681747499
$ws.Range("J8").Value = 1.034990096813273
$ws.Range("K8").Value = 1.041272272745269
$ws.Range("L8").Value = 1.029199422810748
$ws.Range("N8").Value = 1.015540295266384
$ws.Range("Q8").Value = 1.02
$ws.Range("R8").Value = 1.040276209258421
$ws.Range("C9").Value = 1.004437579004822
$ws.Range("D9").Value = 1.024254181176485
$ws.Range("E9").Value = 1.010922759302285
$ws.Range("I9").Value = 1.048483858226924
$ws.Range("J9").Value = 1.029834602248235
$ws.Range("K9").Value = 1.036965636532218
$ws.Range("L9").Value = 1.023841073665792
$ws.Range("N9").Value = 1.013821450940855
$ws.Range("Q9").Value = 1.02
$ws.Range("R9").Value = 1.037227965595801
$ws.Range("C10").Value = 0.9979401178270427
$ws.Range("D10").Value = 1.019939085145463
$ws.Range("E10").Value = 1.005850633227417
$ws.Range("I10").Value = 1.046503709510523
$ws.Range("J10").Value = 1.026235701699336
$ws.Range("K10").Value = 1.033957429666773
$ws.Range("L10").Value = 1.020115050109036
$ws.Range("N10").Value = 1.012623690023414
$ws.Range("Q10").Value = 1.02
$ws.Range("R10").Value = 1.035117782198548
$ws.Range("C11").Value = 0.9950699379343589
$ws.Range("D11").Value = 1.018046955090412
$ws.Range("E11").Value = 1.003621072064988
$ws.Range("I11").Value = 1.045623562315754
$ws.Range("J11").Value = 1.02464926128664
$ws.Range("K11").Value = 1.03263637276962
$ws.Range("L11").Value = 1.018475021837476
$ws.Range("N11").Value = 1.012098888372661
$ws.Range("Q11").Value = 1.02
$ws.Range("R11").Value = 1.03421660122505
$ws.Range("C12").Value = 0.9939844136262626
$ws.Range("D12").Value = 1.017327713307744
$ws.Range("E12").Value = 1.002777910986824
$ws.Range("I12").Value = 1.045285827389805
$ws.Range("J12").Value = 1.024045564040504
$ws.Range("K12").Value = 1.032130304425527
$ws.Range("L12").Value = 1.017852083422442
$ws.Range("N12").Value = 1.011897685321336
$ws.Range("Q12").Value = 1.02
$ws.Range("R12").Value = 1.033858793743822
$ws.Range("C13").Value = 0.9942161744180513
$ws.Range("D13").Value = 1.01748034365777
$ws.Range("E13").Value = 1.002957611105075
$ws.Range("I13").Value = 1.045357388400092
$ws.Range("J13").Value = 1.024173824759202
$ws.Range("K13").Value = 1.032237203173613
$ws.Range("L13").Value = 1.01798452219907
$ws.Range("N13").Value = 1.011940124495672
$ws.Range("Q13").Value = 1.02
$ws.Range("R13").Value = 1.033931878504235
$ws.Range("C14").Value = 0.994979755136962
$ws.Range("D14").Value = 1.017986803882488
$ws.Range("E14").Value = 1.003550888177998
$ws.Range("I14").Value = 1.045595272255369
$ws.Range("J14").Value = 1.024598838753988
$ws.Range("K14").Value = 1.032593839272375
$ws.Range("L14").Value = 1.018423030601334
$ws.Range("N14").Value = 1.012081952128195
$ws.Range("Q14").Value = 1.02
$ws.Range("R14").Value = 1.034185459286751
$ws.Range("C15").Value = 0.9954519004846002
$ws.Range("D15").Value = 1.018301853298875
$ws.Range("E15").Value = 1.003918411785225
$ws.Range("I15").Value = 1.04574338114158
$ws.Range("J15").Value = 1.02486287507806
$ws.Range("K15").Value = 1.032816626427643
$ws.Range("L15").Value = 1.018695290025004
$ws.Range("N15").Value = 1.012170670340454
$ws.Range("Q15").Value = 1.02
$ws.Range("R15").Value = 1.034348856623821
$ws.Range("C16").Value = 0.9981566427875026
$ws.Range("D16").Value = 1.020097363754772
$ws.Range("E16").Value = 1.006023328353933
$ws.Range("I16").Value = 1.046580714672315
$ws.Range("J16").Value = 1.026366568121388
$ws.Range("K16").Value = 1.034077352665794
$ws.Range("L16").Value = 1.020248314026105
$ws.Range("N16").Value = 1.012672193820598
$ws.Range("Q16").Value = 1.02
$ws.Range("R16").Value = 1.035243404807702
$ws.Range("C17").Value = 0.9998335858778766
$ws.Range("D17").Value = 1.021212147895048
$ws.Range("E17").Value = 1.007330692519585
$ws.Range("I17").Value = 1.047096903931525
$ws.Range("J17").Value = 1.027298160637506
$ws.Range("K17").Value = 1.034858199947542
$ws.Range("L17").Value = 1.021211397518679
$ws.Range("N17").Value = 1.012982968276462
$ws.Range("Q17").Value = 1.02
$ws.Range("R17").Value = 1.035798080347772
$ws.Range("C18").Value = 1.000796316469365
$ws.Range("D18").Value = 1.021847957997101
$ws.Range("E18").Value = 1.008080771478984
$ws.Range("I18").Value = 1.047388822430985
$ws.Range("J18").Value = 1.027829273645991
$ws.Range("K18").Value = 1.035299930669959
$ws.Range("L18").Value = 1.021761463014982
$ws.Range("N18").Value = 1.013158608155542
$ws.Range("Q18").Value = 1.02
$ws.Range("R18").Value = 1.036098755093401
$ws.Range("C19").Value = 1.001129981614718
$ws.Range("D19").Value = 1.022072087598616
$ws.Range("E19").Value = 1.008341874864928
$ws.Range("I19").Value = 1.047492423470099
$ws.Range("J19").Value = 1.028015989455178
$ws.Range("K19").Value = 1.035457841363275
$ws.Range("L19").Value = 1.021954374171227
$ws.Range("N19").Value = 1.013221581039357
$ws.Range("Q19").Value = 1.02
$ws.Range("R19").Value = 1.036216798231291
$ws.Range("C20").Value = 0.9996537267542656
$ws.Range("D20").Value = 1.021092159202299
$ws.Range("E20").Value = 1.007190286451166
$ws.Range("I20").Value = 1.047041392906381
$ws.Range("J20").Value = 1.027198003183171
$ws.Range("K20").Value = 1.034773997625085
$ws.Range("L20").Value = 1.021107870753111
$ws.Range("N20").Value = 1.012949432173993
$ws.Range("Q20").Value = 1.02
$ws.Range("R20").Value = 1.035737210275163
$ws.Range("C21").Value = 0.9947642484818813
$ws.Range("D21").Value = 1.017848810226919
$ws.Range("E21").Value = 1.003384829994356
$ws.Range("I21").Value = 1.045531657254999
$ws.Range("J21").Value = 1.024482517149476
$ws.Range("K21").Value = 1.032499741920467
$ws.Range("L21").Value = 1.01830235823846
$ws.Range("N21").Value = 1.012044845845869
$ws.Range("Q21").Value = 1.02
$ws.Range("R21").Value = 1.034132459227193
$ws.Range("C22").Value = 0.9916127080113906
$ws.Range("D22").Value = 1.015757446249093
$ws.Range("E22").Value = 1.000938097994355
$ws.Range("I22").Value = 1.044544514427519
$ws.Range("J22").Value = 1.022725493366126
$ws.Range("K22").Value = 1.031023077412778
$ws.Range("L22").Value = 1.016490993858964
$ws.Range("N22").Value = 1.011457593231728
$ws.Range("Q22").Value = 1.02
$ws.Range("R22").Value = 1.033074874644792
$ws.Range("C23").Value = 0.9932790315541323
$ws.Range("D23").Value = 1.016856750500886
$ws.Range("E23").Value = 1.002229293258094
$ws.Range("I23").Value = 1.045063199836559
$ws.Range("J23").Value = 1.02365035848871
$ws.Range("K23").Value = 1.031796253188069
$ws.Range("L23").Value = 1.01744492554219
$ws.Range("N23").Value = 1.011764637331959
$ws.Range("Q23").Value = 1.02
$ws.Range("R23").Value = 1.033611896705424
$ws.Range("C24").Value = 0.9997171566945828
$ws.Range("D24").Value = 1.021124584879156
$ws.Range("E24").Value = 1.007237087226112
$ws.Range("I24").Value = 1.04705396011135
$ws.Range("J24").Value = 1.02722609878793
$ws.Range("K24").Value = 1.03479061410201
$ws.Range("L24").Value = 1.021138292392223
$ws.Range("N24").Value = 1.012955564200424
$ws.Range("Q24").Value = 1.02
$ws.Range("R24").Value = 1.035721548876428
$ws.Range("C25").Value = 1.006917101785286
$ws.Range("D25").Value = 1.025917176138622
$ws.Range("E25").Value = 1.012868631681859
$ws.Range("I25").Value = 1.049237705054562
$ws.Range("J25").Value = 1.031213899371863
$ws.Range("K25").Value = 1.03812605543696
$ws.Range("L25").Value = 1.02527041003154
$ws.Range("N25").Value = 1.014284405465614
$ws.Range("Q25").Value = 1.02
$ws.Range("R25").Value = 1.038077012395498
